# prepare for 2024qld live forecast
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the scenario/forecast label from "qld" to "fed" (M5),
# which flips the IF($M$5="qld", ...) formulas in row 7 to their
# federal (non-qld) branch values.
$ws.Range("M5").Value = "fed"

# New first-preference poll numbers in row 2.
$ws.Range("A2").Value = 35
$ws.Range("B2").Value = 28
$ws.Range("C2").Value = 12
$ws.Range("D2").Value = 7
$ws.Range("F2").Value = 2

# I2 used to be "100 - (the other figures)"; it's now a plain entered
# figure (the formula is replaced with a literal value).
$ws.Range("I2").Value = 9

# Second set of figures (2PP-ish inputs) in row 10.
$ws.Range("A10").Value = 48
$ws.Range("B10").Value = 46

# Re-point the active selection to A4:I4 (matches the new sheetView state).
[void]$ws.Range("A4:I4").Select()
